# Evidencias workbook update: real data + new "Projeto" lookup columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new columns needed -------------------------------------
# New column for "Nome do Ponto" right after "Id do Ponto" (old column H)
$ws.Columns("I:I").Insert()
# Two new columns for "Id do Projeto" / "Nome do Projeto" right before the
# (now shifted) "Usuário que criou" column (which sits at M after the first insert)
$ws.Columns("M:N").Insert()

# --- 2. Header row (row 1) --------------------------------------------------
$ws.Range("A1").Value = "Id da Evidência"
$ws.Range("B1").Value = "Status da Evidência"
$ws.Range("C1").Value = "Tipo"
$ws.Range("D1").Value = "Quantidade"
$ws.Range("E1").Value = "Notas"
$ws.Range("F1").Value = "Profundidade"
$ws.Range("G1").Value = "Solo"
$ws.Range("H1").Value = "Id do Ponto"
$ws.Range("I1").Value = "Nome do Ponto"
$ws.Range("J1").Value = "Longitude do Ponto"
$ws.Range("K1").Value = "Latitude do Ponto"
$ws.Range("L1").Value = "Status do Ponto"
$ws.Range("M1").Value = "Id do Projeto"
$ws.Range("N1").Value = "Nome do Projeto"
$ws.Range("O1").Value = "Usuário que criou"
$ws.Range("P1").Value = "Usuário que atualizou"
$ws.Range("Q1").Value = "Data de criação"
$ws.Range("R1").Value = "Data de atualização"

# --- 3. Fill in the new "Nome do Ponto" values (I2:I5) ----------------------
# These look like plain numbers but are stored as text in the source system,
# so the column is pre-formatted as Text before the values are typed in.
$ws.Range("I2:I5").NumberFormat = "@"
$ws.Range("I2").Value = "139"
$ws.Range("I3").Value = "146"
$ws.Range("I4").Value = "149"
$ws.Range("I5").Value = "151"

# --- 4. Fill in "Id do Ponto" values (H2:H5) --------------------------------
$ws.Range("H2").Value = """5dc21aa103e32600176a26e6"""
$ws.Range("H3").Value = """5dc21aa103e32600176a26ed"""
$ws.Range("H4").Value = """5dc21aa103e32600176a26f0"""
$ws.Range("H5").Value = """5dc21aa103e32600176a26f2"""

# --- 5. Fill in the new "Id do Projeto" / "Nome do Projeto" columns --------
$ws.Range("M2:M5").Value = """5dbf916a598a81001721843c"""
$ws.Range("N2:N5").Value = "Teste"

# --- 6. Update "Usuário que criou" / "Usuário que atualizou" with real names
$ws.Range("O2:O5").Value = "Rodrigo Mota"
$ws.Range("P2").Value = "Rodrigo Mota"
$ws.Range("P3").Value = "Rodrigo Mota"
$ws.Range("P4").Value = "Rodrigo Mota"
$ws.Range("P5").Value = "Gustavo"

# --- 7. Column widths --------------------------------------------------------
# ColumnWidth has a constant +5/6 character padding baked in by the engine,
# so we subtract it to land exactly on the target stored width.
$pad = 0.8333333333333334
$ws.Columns("A:C").ColumnWidth = 30 - $pad
$ws.Columns("D:D").ColumnWidth = 15 - $pad
$ws.Columns("E:E").ColumnWidth = 80 - $pad
$ws.Columns("F:H").ColumnWidth = 30 - $pad
$ws.Columns("I:I").ColumnWidth = 50 - $pad
$ws.Columns("J:J").ColumnWidth = 320 - $pad
$ws.Columns("K:R").ColumnWidth = 30 - $pad
